# Black Scholes with setting parameters in Variables environment
$wb = $excel.ActiveWorkbook

# Rename "Scenario Generator" sheet to "Scenario_Generator"
$scenarioSheet = $wb.Worksheets.Item("Scenario Generator")
$scenarioSheet.Name = "Scenario_Generator"

# Update selection on Input_3M sheet (stays at G16, tab no longer selected)
$input3M = $wb.Worksheets.Item("Input_3M")
$input3M.Range("G16").Select()

# Move selection on Scenario_Generator sheet from J24 to I32, and make it the active tab
$scenarioSheet.Activate()
$scenarioSheet.Range("I32").Select()
